# documentation for candidateGenes change: Level->Confidence score
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row text updates
$ws.Range("A1").Value = "Gene"
$ws.Range("B1").Value = "pathologyID"
$ws.Range("C1").Value = "Confidence score"

# New explicit column widths for columns B and C (column A keeps its
# existing explicit width of 19.31 characters, unchanged by this edit).
$ws.Columns.Item(2).ColumnWidth = 13.0
$ws.Columns.Item(3).ColumnWidth = 19.666666666666668
